$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (previously rows 8-10, "Resolving-Mac" target cluster)
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bmp6"
$ws.Cells.Item(2, 3).Value = "Bmpr1b"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 5.922420666666667
$ws.Cells.Item(2, 8).Value = 17.767262
$ws.Cells.Item(2, 9).Value = 0.5833698282960434
$ws.Cells.Item(2, 10).Value = 0.6311054116979437
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.467406
$ws.Cells.Item(2, 14).Value = 4.402218
$ws.Cells.Item(2, 15).Value = 0.864087546066766
$ws.Cells.Item(2, 16).Value = 0.9050919696083439
$ws.Cells.Item(2, 17).Value = 8.690595620790665
$ws.Cells.Item(2, 18).Value = 78.21536058711598
$ws.Cells.Item(2, 19).Value = 0.5040826033817188
$ws.Cells.Item(2, 20).Value = 0.5712084401041766

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bmp6"
$ws.Cells.Item(3, 3).Value = "Bmpr1b"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 5.922420666666667
$ws.Cells.Item(3, 8).Value = 17.767262
$ws.Cells.Item(3, 9).Value = 0.5833698282960434
$ws.Cells.Item(3, 10).Value = 0.6311054116979437
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.2308085
$ws.Cells.Item(3, 14).Value = 0.461617
$ws.Cells.Item(3, 15).Value = 0.135912453933234
$ws.Cells.Item(3, 16).Value = 0.09490803039165596
$ws.Cells.Item(3, 17).Value = 1.366945030442333
$ws.Cells.Item(3, 18).Value = 8.201670182653999
$ws.Cells.Item(3, 19).Value = 0.07928722491432463
$ws.Cells.Item(3, 20).Value = 0.05989697159376699

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Bmp6"
$ws.Cells.Item(4, 3).Value = "Bmpr1b"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.926013333333334
$ws.Cells.Item(4, 8).Value = 5.778040000000001
$ws.Cells.Item(4, 9).Value = 0.1897160182974547
$ws.Cells.Item(4, 10).Value = 0.2052399696141807
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.467406
$ws.Cells.Item(4, 14).Value = 4.402218
$ws.Cells.Item(4, 15).Value = 0.864087546066766
$ws.Cells.Item(4, 16).Value = 0.9050919696083439
$ws.Cells.Item(4, 17).Value = 2.826243521413333
$ws.Cells.Item(4, 18).Value = 25.43619169272
$ws.Cells.Item(4, 19).Value = 0.1639312487002053
$ws.Cells.Item(4, 20).Value = 0.1857610483404555

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Bmp6"
$ws.Cells.Item(5, 3).Value = "Bmpr1b"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.926013333333334
$ws.Cells.Item(5, 8).Value = 5.778040000000001
$ws.Cells.Item(5, 9).Value = 0.1897160182974547
$ws.Cells.Item(5, 10).Value = 0.2052399696141807
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.2308085
$ws.Cells.Item(5, 14).Value = 0.461617
$ws.Cells.Item(5, 15).Value = 0.135912453933234
$ws.Cells.Item(5, 16).Value = 0.09490803039165596
$ws.Cells.Item(5, 17).Value = 0.4445402484466667
$ws.Cells.Item(5, 18).Value = 2.66724149068
$ws.Cells.Item(5, 19).Value = 0.02578476959724939
$ws.Cells.Item(5, 20).Value = 0.01947892127372521

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Bmp6"
$ws.Cells.Item(6, 3).Value = "Bmpr1b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.303652
$ws.Cells.Item(6, 8).Value = 4.607303999999999
$ws.Cells.Item(6, 9).Value = 0.2269141534065018
$ws.Cells.Item(6, 10).Value = 0.1636546186878756
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.467406
$ws.Cells.Item(6, 14).Value = 4.402218
$ws.Cells.Item(6, 15).Value = 0.864087546066766
$ws.Cells.Item(6, 16).Value = 0.9050919696083439
$ws.Cells.Item(6, 17).Value = 3.380392766711999
$ws.Cells.Item(6, 18).Value = 20.282356600272
$ws.Cells.Item(6, 19).Value = 0.1960736939848418
$ws.Cells.Item(6, 20).Value = 0.1481224811637118

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Bmp6"
$ws.Cells.Item(7, 3).Value = "Bmpr1b"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.303652
$ws.Cells.Item(7, 8).Value = 4.607303999999999
$ws.Cells.Item(7, 9).Value = 0.2269141534065018
$ws.Cells.Item(7, 10).Value = 0.1636546186878756
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.2308085
$ws.Cells.Item(7, 14).Value = 0.461617
$ws.Cells.Item(7, 15).Value = 0.135912453933234
$ws.Cells.Item(7, 16).Value = 0.09490803039165596
$ws.Cells.Item(7, 17).Value = 0.5317024626419999
$ws.Cells.Item(7, 18).Value = 2.126809850568
$ws.Cells.Item(7, 19).Value = 0.03084045942165997
$ws.Cells.Item(7, 20).Value = 0.01553213752416377
